$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 3.182878228561681
$ws.Range("C2").Value = 1.65323645889881
$ws.Range("D2").Value = 157.8057217802531
$ws.Range("E2").Value = 71517.89157740913
$ws.Range("G2").Value = 71680.53341387684

$ws.Range("B3").Value = 1.505614041169197
$ws.Range("C3").Value = 1.65323645889881
$ws.Range("D3").Value = 0.7127328510149897
$ws.Range("E3").Value = 6.48142807727062
$ws.Range("G3").Value = 10.35301142835362

$ws.Range("B4").Value = 0.7287194209349384
$ws.Range("C4").Value = 1.65323645889881
$ws.Range("D4").Value = 0.1529057820181812
$ws.Range("E4").Value = 0.4998867070740569
$ws.Range("G4").Value = 3.034748368925986

$ws.Range("B5").Value = 1.505614041169197
$ws.Range("C5").Value = 1.65323645889881
$ws.Range("D5").Value = 0.1529057820181812
$ws.Range("E5").Value = 0.4998867070740569
$ws.Range("G5").Value = 3.811642989160245

$ws.Range("B6").Value = 0.02258322285507441
$ws.Range("C6").Value = 0.004309184025731883
$ws.Range("D6").Value = 157.8057217802531
$ws.Range("E6").Value = 246.9852506941017
$ws.Range("G6").Value = 404.8178648812356

$ws.Range("B7").Value = 0.006876353814593728
$ws.Range("C7").Value = 0.000002220651329265522
$ws.Range("D7").Value = 2938.103010863317
$ws.Range("E7").Value = 246.9852506941017
$ws.Range("G7").Value = 3185.095140131885

$ws.Range("B8").Value = 0.7287194209349384
$ws.Range("C8").Value = 1.65323645889881
$ws.Range("D8").Value = 157.8057217802531
$ws.Range("E8").Value = 71517.89157740913
$ws.Range("G8").Value = 71678.07925506921
